# Commit after jenkins integration
#
# The workbook originally has a single sheet "Sheet1" holding a small user
# table (UserID/UserName/FirstName/LastName/Email/Password/Phone) in rows
# 1-4. This edit:
#   1. Splits the old rows 3 & 4 (testuser2/testuser3) off onto a brand new
#      "Sheet2" (placed after Sheet1), keeping Sheet1's row 2 untouched.
#   2. Replaces Sheet1 rows 3-4 with new user data and appends a new row 5.
#   3. Turns the Email/Password cells of the new Sheet1 rows into real
#      mailto: hyperlinks (Email & Password both contain "@", so Excel's
#      "autoformat as you type" would have linkified both).
#   4. Resizes columns E/F/G on Sheet1 and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New Sheet2, inserted right after Sheet1, carrying off the old
#    testuser2 / testuser3 rows (with their original formatting).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)

# Copy header formatting (row 1) and data-row formatting (rows 3 & 4,
# which still hold the pre-edit values at this point) from Sheet1 so
# Sheet2 ends up with the same cell styles without inventing new ones.
$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)
$ws1.Range("A3:G4").Copy()
$ws2.Range("A2:G3").PasteSpecial(-4122)
$ws2.Range("A1").Select()
$excel.CutCopyMode = 0

$ws2.Range("A1").Value = "UserID"
$ws2.Range("B1").Value = "UserName"
$ws2.Range("C1").Value = "FirstName"
$ws2.Range("D1").Value = "LastName"
$ws2.Range("E1").Value = "Email"
$ws2.Range("F1").Value = "Password"
$ws2.Range("G1").Value = "Phone"

$ws2.Range("A2").Value = "'1020"
$ws2.Range("B2").Value = "testuser2"
$ws2.Range("C2").Value = "Mohan"
$ws2.Range("D2").Value = "Sharma"
$ws2.Range("E2").Value = "kirankatkar3318@gmail.com"
$ws2.Range("F2").Value = "Kiran@3318"
$ws2.Range("G2").Value = 1234567893

$ws2.Range("A3").Value = "'1030"
$ws2.Range("B3").Value = "testuser3"
$ws2.Range("C3").Value = "Sohan"
$ws2.Range("D3").Value = "Sharma"
$ws2.Range("E3").Value = "kirankatkar3318@gmail.com"
$ws2.Range("F3").Value = "Kiran@3318"
$ws2.Range("G3").Value = 1234567893

# ---------------------------------------------------------------------
# 2. Sheet1: replace rows 3 & 4, add new row 5.
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = "'1011"
$ws1.Range("B3").Value = "testuser2"
$ws1.Range("C3").Value = "aniket"
$ws1.Range("D3").Value = "varma"
$ws1.Range("E3").Value = "sakiran2230@gmail.com"
$ws1.Range("F3").Value = "Kiran@2230"
$ws1.Range("G3").Value = 21343113134

$ws1.Range("A4").Value = "'1013"
$ws1.Range("B4").Value = "testuser3"
$ws1.Range("C4").Value = "abhi"
$ws1.Range("D4").Value = "katkar"
$ws1.Range("E4").Value = "kiran.akatkar3317@gmail.com"
$ws1.Range("F4").Value = "Kiran@9103"
$ws1.Range("G4").Value = 1234567893

$ws1.Range("A5").Value = "'1014"
$ws1.Range("B5").Value = "testusr"
$ws1.Range("C5").Value = "kiran"
$ws1.Range("D5").Value = "katkar"
$ws1.Range("E5").Value = "kirankatkar3318@gmail.com"
$ws1.Range("F5").Value = "Kiran@3318"
$ws1.Range("G5").Value = 7328832329

# ---------------------------------------------------------------------
# 3. Real hyperlinks for the new Email / Password cells, re-applying the
#    "Hyperlink" cell style afterwards so every linked cell keeps using
#    the same style index the pre-existing E2/F2 pair already used.
# ---------------------------------------------------------------------
$linkCells = @("E3","F3","E4","F4","E5","F5")
foreach ($addr in $linkCells) {
    $cell = $ws1.Range($addr)
    $ws1.Hyperlinks.Add($cell, "mailto:" + $cell.Value)
}
foreach ($addr in $linkCells) {
    $ws1.Range($addr).Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 4. Column widths on Sheet1 (E/F/G) and the active selection.
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 20.6666666
$ws1.Columns.Item(6).ColumnWidth = 16.1666666
$ws1.Columns.Item(7).ColumnWidth = 11.1666666

$ws1.Range("F3").Select()
$ws1.Activate()
